$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 770.1667
$ws.Range("I15").Value = 770.1667
$ws.Range("K15").Value = 2310.5001
$ws.Range("M15").Value = -2141.5001
$ws.Range("H74").Value = 6388.3335
$ws.Range("I74").Value = 6388.3335
$ws.Range("K74").Value = 6388.3335
$ws.Range("M74").Value = -5452.3335
$ws.Range("H76").Value = 5095.4
$ws.Range("I76").Value = 3200
$ws.Range("K76").Value = 3200
$ws.Range("M76").Value = -2885
$ws.Range("H77").Value = 6388.3335
$ws.Range("I77").Value = 6388.3335
$ws.Range("K77").Value = 31941.6675
$ws.Range("M77").Value = -27261.6675
$ws.Range("H79").Value = 5095.4
$ws.Range("I79").Value = 3200
$ws.Range("K79").Value = 3200
$ws.Range("M79").Value = -2108
$ws.Range("H92").Value = 493.2143
$ws.Range("I92").Value = 463.07693
$ws.Range("K92").Value = 463.07693
$ws.Range("M92").Value = 784.9230700000001
$ws.Range("H112").Value = 2261.5227
$ws.Range("J112").Value = 2318.95
$ws.Range("L112").Value = 6956.849999999999
$ws.Range("N112").Value = -9172.849999999999
$ws.Range("H131").Value = 2886
$ws.Range("I131").Value = 2886
$ws.Range("K131").Value = 8658
$ws.Range("M131").Value = -3618
$ws.Range("H133").Value = 59388.332
$ws.Range("J133").Value = 59388.332
$ws.Range("L133").Value = 59388.332
$ws.Range("N133").Value = -69508.33199999999
$ws.Range("H137").Value = 7833.4067
$ws.Range("I137").Value = 6446.7
$ws.Range("J137").Value = 9267.931
$ws.Range("K137").Value = 19340.1
$ws.Range("L137").Value = 27803.793
$ws.Range("M137").Value = -16790.1
$ws.Range("N137").Value = -32903.79300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 27547.684
$ws.Range("I97").Value = 28251.676
$ws.Range("K97").Value = 28251.676
$ws.Range("M97").Value = -27755.676
$ws.Range("H122").Value = 3087.5
$ws.Range("I122").Value = 3107
$ws.Range("K122").Value = 9321
$ws.Range("M122").Value = -6871

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 4335.25
$ws.Range("I10").Value = 4294
$ws.Range("J10").Value = 4349
$ws.Range("K10").Value = 4294
$ws.Range("L10").Value = 4349
$ws.Range("M10").Value = -4154
$ws.Range("N10").Value = -4629
$ws.Range("H74").Value = 59932
$ws.Range("J74").Value = 59932
$ws.Range("L74").Value = 59932
$ws.Range("N74").Value = -61804
$ws.Range("H77").Value = 59932
$ws.Range("J77").Value = 59932
$ws.Range("L77").Value = 179796
$ws.Range("N77").Value = -189156
$ws.Range("H99").Value = 2016
$ws.Range("I99").Value = 2016
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2016
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -518
$ws.Range("H132").Value = 79760
$ws.Range("J132").Value = 79760
$ws.Range("L132").Value = 79760
$ws.Range("N132").Value = -89880
$ws.Range("H134").Value = 6971.625
$ws.Range("J134").Value = 17649.334
$ws.Range("L134").Value = 52948.00199999999
$ws.Range("N134").Value = -58018.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12548.154
$ws.Range("J58").Value = 20680
$ws.Range("L58").Value = 20680
$ws.Range("N58").Value = -21086
$ws.Range("H94").Value = 1107.9412
$ws.Range("I94").Value = 1146.4
$ws.Range("K94").Value = 1146.4
$ws.Range("M94").Value = -695.4000000000001
$ws.Range("H132").Value = 22179.09
$ws.Range("I132").Value = 16625.625
$ws.Range("J132").Value = 29583.709
$ws.Range("K132").Value = 49876.875
$ws.Range("L132").Value = 88751.12699999999
$ws.Range("M132").Value = -47346.875
$ws.Range("N132").Value = -93811.12699999999
$ws.Range("H134").Value = 8837.659
$ws.Range("J134").Value = 16135.546
$ws.Range("L134").Value = 48406.638
$ws.Range("N134").Value = -53476.638
$ws.Range("H136").Value = 12548.154
$ws.Range("J136").Value = 20680
$ws.Range("L136").Value = 62040
$ws.Range("N136").Value = -67140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1796.7142
$ws.Range("J5").Value = 2600.3157
$ws.Range("L5").Value = 7800.9471
$ws.Range("N5").Value = -8024.9471
$ws.Range("H14").Value = 4498.3
$ws.Range("I14").Value = 4498.3
$ws.Range("K14").Value = 13494.9
$ws.Range("M14").Value = -13321.9
$ws.Range("H56").Value = 3971.3333
$ws.Range("I56").Value = 3971.3333
$ws.Range("K56").Value = 3971.3333
$ws.Range("M56").Value = -3441.3333
$ws.Range("H86").Value = 797.58826
$ws.Range("I86").Value = 791.8333
$ws.Range("J86").Value = 800.7273
$ws.Range("K86").Value = 2375.4999
$ws.Range("L86").Value = 2402.1819
$ws.Range("M86").Value = -1189.4999
$ws.Range("N86").Value = -4774.1819
$ws.Range("H89").Value = 797.58826
$ws.Range("I89").Value = 791.8333
$ws.Range("J89").Value = 800.7273
$ws.Range("K89").Value = 7126.4997
$ws.Range("L89").Value = 7206.545700000001
$ws.Range("M89").Value = -1198.4997
$ws.Range("N89").Value = -19062.5457
$ws.Range("H107").Value = 1391.091
$ws.Range("J107").Value = 1362.2667
$ws.Range("L107").Value = 4086.800099999999
$ws.Range("N107").Value = -7926.800099999999
$ws.Range("H119").Value = 5999
$ws.Range("I119").Value = 5999
$ws.Range("K119").Value = 17997
$ws.Range("M119").Value = -13159
$ws.Range("H135").Value = 1796.7142
$ws.Range("J135").Value = 2600.3157
$ws.Range("L135").Value = 23402.8413
$ws.Range("N135").Value = -28472.8413
$ws.Range("H138").Value = 3400
$ws.Range("I138").Value = 3400
$ws.Range("K138").Value = 10200
$ws.Range("M138").Value = -5060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14498.8125
$ws.Range("J132").Value = 31198.8
$ws.Range("L132").Value = 93596.39999999999
$ws.Range("N132").Value = -98656.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1145.8928
$ws.Range("I22").Value = 1371.4706
$ws.Range("J22").Value = 797.2727
$ws.Range("K22").Value = 1371.4706
$ws.Range("L22").Value = 797.2727
$ws.Range("M22").Value = -1076.4706
$ws.Range("N22").Value = -1387.2727
$ws.Range("H27").Value = 1145.8928
$ws.Range("I27").Value = 1371.4706
$ws.Range("J27").Value = 797.2727
$ws.Range("K27").Value = 1371.4706
$ws.Range("L27").Value = 797.2727
$ws.Range("M27").Value = -1264.4706
$ws.Range("N27").Value = -1011.2727
$ws.Range("H68").Value = 3703.913
$ws.Range("I68").Value = 2694.1177
$ws.Range("J68").Value = 6565
$ws.Range("K68").Value = 2694.1177
$ws.Range("L68").Value = 6565
$ws.Range("M68").Value = -1945.1177
$ws.Range("N68").Value = -8063
$ws.Range("H71").Value = 3703.913
$ws.Range("I71").Value = 2694.1177
$ws.Range("J71").Value = 6565
$ws.Range("K71").Value = 13470.5885
$ws.Range("L71").Value = 32825
$ws.Range("M71").Value = -9726.588499999998
$ws.Range("N71").Value = -40313
$ws.Range("H132").Value = 2855733.5
$ws.Range("I132").Value = 3708821.2
$ws.Range("J132").Value = 12108
$ws.Range("K132").Value = 11126463.6
$ws.Range("L132").Value = 36324
$ws.Range("M132").Value = -11123933.6
$ws.Range("N132").Value = -41384
$ws.Range("H136").Value = 4448507.5
$ws.Range("I136").Value = 5559109
$ws.Range("J136").Value = 6100
$ws.Range("K136").Value = 16677327
$ws.Range("L136").Value = 18300
$ws.Range("M136").Value = -16674777
$ws.Range("N136").Value = -23400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 9679.333000000001
$ws.Range("J31").Value = 10019
$ws.Range("L31").Value = 10019
$ws.Range("N31").Value = -10715
$ws.Range("H97").Value = 35572
$ws.Range("J97").Value = 35572
$ws.Range("L97").Value = 35572
$ws.Range("N97").Value = -37554
$ws.Range("H126").Value = 63834.332
$ws.Range("I126").Value = 8876.1875
$ws.Range("K126").Value = 26628.5625
$ws.Range("M126").Value = -24158.5625
$ws.Range("H130").Value = 96982
$ws.Range("J130").Value = 96982
$ws.Range("L130").Value = 96982
$ws.Range("N130").Value = -107022
$ws.Range("H132").Value = 18569.94
$ws.Range("I132").Value = 16501.727
$ws.Range("K132").Value = 49505.181
$ws.Range("M132").Value = -46975.181
$ws.Range("H141").Value = 63747.5
$ws.Range("J141").Value = 63747.5
$ws.Range("L141").Value = 63747.5
$ws.Range("N141").Value = -74107.5
